$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 38.55267666666666
$ws.Range("N2").Value = 115.65803
$ws.Range("O2").Value = 0.5758151725879548
$ws.Range("P2").Value = 0.5758151725879548
$ws.Range("Q2").Value = 1.322485318588889
$ws.Range("R2").Value = 11.9023678673
$ws.Range("S2").Value = 0.5758151725879548
$ws.Range("T2").Value = 0.5758151725879548

# Row 3
$ws.Range("O3").Value = 0.08021535714867321
$ws.Range("P3").Value = 0.08021535714867323
$ws.Range("S3").Value = 0.08021535714867321
$ws.Range("T3").Value = 0.08021535714867323

# Row 4
$ws.Range("M4").Value = 23.02986166666667
$ws.Range("N4").Value = 69.089585
$ws.Range("O4").Value = 0.3439694702633719
$ws.Range("P4").Value = 0.3439694702633719
$ws.Range("Q4").Value = 0.7900010213722222
$ws.Range("R4").Value = 7.11000919235
$ws.Range("S4").Value = 0.3439694702633719
$ws.Range("T4").Value = 0.3439694702633719
